$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New locale row (en-ES) -------------------------------------------------
# Fill the WelcomeMessage..ActiveDropInScreendescription columns first (C:G),
# then the Locale (A) and IsSupported (B) columns, matching the order the
# strings were originally authored in.
$ws.Range("C3").Value = "`nBonjour,"
$ws.Range("C3").WrapText = $true
$ws.Range("D3").Value = "Nouveau drop-in"
$ws.Range("E3").Value = "Visiter une classe virtuelle"
$ws.Range("F3").Value = "Drop-ins programmés"
$ws.Range("G3").Value = "Afficher les détails ou prolonger la durée"
$ws.Range("A3").Value = "en-ES"
$ws.Range("B3").Value = $true

$ws.Rows.Item(3).RowHeight = 28.8

# --- New columns: owner/member drop-in labels -------------------------------
$tbl = $ws.ListObjects.Item(1)

$tbl.ListColumns.Add() | Out-Null
$ws.Range("CL1").Value = "OwnerOrMemberLbl"
$ws.Range("CL2").Value = "Dropping the user as :"

$tbl.ListColumns.Add() | Out-Null
$ws.Range("CM1").Value = "OwnerLbl"
$ws.Range("CM2").Value = "Owner"

$tbl.ListColumns.Add() | Out-Null
$ws.Range("CN1").Value = "MemberLbl"
$ws.Range("CN2").Value = "Member"

# --- Selection reflects the newly edited area -------------------------------
$ws.Range("CI19").Select() | Out-Null
